$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts everything below down by one)
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0.1
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = 0.9
$ws.Range("G8").Value = 0.4
$ws.Range("H8").Value = 0.4

# Insert a new row at position 13 (before franzosa_ControlvsUC_Fp), shifts everything below down by one
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "franzosa_ControlvsUC_ConvUC"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.6
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4
$ws.Range("H13").Value = 0.4

Write-Output "Done"
